$d = $word.ActiveDocument

# The first paragraph currently reads:
#   "This is a Microsoft word document."
# We need to append " (Changed main)" as three additional, separate
# runs (matching the target OOXML run layout) after the existing run,
# without touching the original run's text or formatting.

$firstPara = $d.Paragraphs(1).Range
# Range covering just the paragraph's text (exclude the paragraph mark).
$textRange = $d.Range($firstPara.Start, $firstPara.End - 1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p>' +
       '<w:r><w:t>This is a Microsoft word document.</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
       '<w:r><w:t>Changed main</w:t></w:r>' +
       '<w:r><w:t>)</w:t></w:r>' +
       '</w:p></w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$textRange.InsertXML($xml)
